$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(97, 8).Value = 3466.6667
$ws.Cells.Item(97, 10).Value = 3466.6667
$ws.Cells.Item(97, 12).Value = 10400.0001
$ws.Cells.Item(97, 14).Value = -11392.0001
$ws.Cells.Item(98, 8).Value = 1763
$ws.Cells.Item(98, 9).Value = 1249.6154
$ws.Cells.Item(98, 11).Value = 1249.6154
$ws.Cells.Item(98, 13).Value = 248.3846000000001
$ws.Cells.Item(116, 8).Value = 4678.9565
$ws.Cells.Item(116, 9).Value = 4362.7144
$ws.Cells.Item(116, 11).Value = 4362.7144
$ws.Cells.Item(116, 13).Value = -920.7143999999998
$ws.Cells.Item(122, 8).Value = 1763
$ws.Cells.Item(122, 9).Value = 1249.6154
$ws.Cells.Item(122, 11).Value = 3748.8462
$ws.Cells.Item(122, 13).Value = -1298.8462
$ws.Cells.Item(125, 8).Value = 11839011
$ws.Cells.Item(125, 9).Value = 3390947.8
$ws.Cells.Item(125, 10).Value = 17873342
$ws.Cells.Item(125, 11).Value = 30518530.2
$ws.Cells.Item(125, 12).Value = 160860078
$ws.Cells.Item(125, 13).Value = -30516070.2
$ws.Cells.Item(125, 14).Value = -160864998
$ws.Cells.Item(129, 8).Value = 10123.5
$ws.Cells.Item(129, 9).Value = 1693.1111
$ws.Cells.Item(129, 10).Value = 85997
$ws.Cells.Item(129, 11).Value = 5079.3333
$ws.Cells.Item(129, 12).Value = 257991
$ws.Cells.Item(129, 13).Value = -79.33330000000024
$ws.Cells.Item(129, 14).Value = -267991
$ws.Cells.Item(137, 8).Value = 1418.4
$ws.Cells.Item(137, 9).Value = 1328.8572
$ws.Cells.Item(137, 10).Value = 1627.3334
$ws.Cells.Item(137, 11).Value = 3986.5716
$ws.Cells.Item(137, 12).Value = 4882.0002
$ws.Cells.Item(137, 13).Value = -1436.5716
$ws.Cells.Item(137, 14).Value = -9982.0002
$ws.Cells.Item(138, 8).Value = 4643.88
$ws.Cells.Item(138, 10).Value = 5427.923
$ws.Cells.Item(138, 12).Value = 16283.769
$ws.Cells.Item(138, 14).Value = -26563.769
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4522.5
$ws.Cells.Item(32, 9).Value = 2023.4054
$ws.Cells.Item(32, 10).Value = 23015.8
$ws.Cells.Item(32, 11).Value = 2023.4054
$ws.Cells.Item(32, 12).Value = 23015.8
$ws.Cells.Item(32, 13).Value = -1736.4054
$ws.Cells.Item(32, 14).Value = -23589.8
$ws.Cells.Item(61, 8).Value = 33335202
$ws.Cells.Item(61, 9).Value = 35715896
$ws.Cells.Item(61, 10).Value = 5514.5
$ws.Cells.Item(61, 11).Value = 35715896
$ws.Cells.Item(61, 12).Value = 5514.5
$ws.Cells.Item(61, 13).Value = -35715684
$ws.Cells.Item(61, 14).Value = -5938.5
$ws.Cells.Item(102, 8).Value = 1632.4
$ws.Cells.Item(102, 9).Value = 1416.75
$ws.Cells.Item(102, 10).Value = 2495
$ws.Cells.Item(102, 11).Value = 1416.75
$ws.Cells.Item(102, 12).Value = 2495
$ws.Cells.Item(102, 13).Value = 205.25
$ws.Cells.Item(102, 14).Value = -5739
$ws.Cells.Item(132, 8).Value = 2440299.8
$ws.Cells.Item(132, 9).Value = 2565390.2
$ws.Cells.Item(132, 10).Value = 1035
$ws.Cells.Item(132, 11).Value = 7696170.600000001
$ws.Cells.Item(132, 12).Value = 3105
$ws.Cells.Item(132, 13).Value = -7693640.600000001
$ws.Cells.Item(132, 14).Value = -8165
$ws.Cells.Item(136, 8).Value = 33335202
$ws.Cells.Item(136, 9).Value = 35715896
$ws.Cells.Item(136, 10).Value = 5514.5
$ws.Cells.Item(136, 11).Value = 107147688
$ws.Cells.Item(136, 12).Value = 16543.5
$ws.Cells.Item(136, 13).Value = -107145138
$ws.Cells.Item(136, 14).Value = -21643.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1848
$ws.Cells.Item(99, 9).Value = 1848
$ws.Cells.Item(99, 11).Value = 1848
$ws.Cells.Item(99, 13).Value = -350
$ws.Cells.Item(106, 8).Value = 70000
$ws.Cells.Item(106, 10).Value = 70000
$ws.Cells.Item(106, 12).Value = 70000
$ws.Cells.Item(106, 14).Value = -72524
$ws.Cells.Item(109, 8).Value = 41666.5
$ws.Cells.Item(109, 9).Value = 25000
$ws.Cells.Item(109, 10).Value = 44999.8
$ws.Cells.Item(109, 11).Value = 25000
$ws.Cells.Item(109, 12).Value = 44999.8
$ws.Cells.Item(109, 13).Value = -23613
$ws.Cells.Item(109, 14).Value = -47773.8
$ws.Cells.Item(134, 8).Value = 11861332
$ws.Cells.Item(134, 9).Value = 12439870
$ws.Cells.Item(134, 10).Value = 1300
$ws.Cells.Item(134, 11).Value = 37319610
$ws.Cells.Item(134, 12).Value = 3900
$ws.Cells.Item(134, 13).Value = -37317075
$ws.Cells.Item(134, 14).Value = -8970
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 8684.833000000001
$ws.Cells.Item(7, 9).Value = 11479.777
$ws.Cells.Item(7, 10).Value = 300
$ws.Cells.Item(7, 11).Value = 11479.777
$ws.Cells.Item(7, 12).Value = 300
$ws.Cells.Item(7, 13).Value = -11366.777
$ws.Cells.Item(7, 14).Value = -526
$ws.Cells.Item(16, 8).Value = 989541.8
$ws.Cells.Item(16, 9).Value = 1088296
$ws.Cells.Item(16, 11).Value = 1088296
$ws.Cells.Item(16, 13).Value = -1088009
$ws.Cells.Item(31, 8).Value = 2970.7173
$ws.Cells.Item(31, 9).Value = 4951.6113
$ws.Cells.Item(31, 10).Value = 1697.2858
$ws.Cells.Item(31, 11).Value = 4951.6113
$ws.Cells.Item(31, 12).Value = 1697.2858
$ws.Cells.Item(31, 13).Value = -4656.6113
$ws.Cells.Item(31, 14).Value = -2287.2858
$ws.Cells.Item(34, 8).Value = 2970.7173
$ws.Cells.Item(34, 9).Value = 4951.6113
$ws.Cells.Item(34, 10).Value = 1697.2858
$ws.Cells.Item(34, 11).Value = 4951.6113
$ws.Cells.Item(34, 12).Value = 1697.2858
$ws.Cells.Item(34, 13).Value = -4749.6113
$ws.Cells.Item(34, 14).Value = -2101.2858
$ws.Cells.Item(86, 8).Value = 9897.625
$ws.Cells.Item(86, 9).Value = 8244.5
$ws.Cells.Item(86, 11).Value = 8244.5
$ws.Cells.Item(86, 13).Value = -7121.5
$ws.Cells.Item(89, 8).Value = 9897.625
$ws.Cells.Item(89, 9).Value = 8244.5
$ws.Cells.Item(89, 11).Value = 41222.5
$ws.Cells.Item(89, 13).Value = -35606.5
$ws.Cells.Item(99, 8).Value = 11509.417
$ws.Cells.Item(99, 9).Value = 13912.667
$ws.Cells.Item(99, 10).Value = 4299.6665
$ws.Cells.Item(99, 11).Value = 13912.667
$ws.Cells.Item(99, 12).Value = 4299.6665
$ws.Cells.Item(99, 13).Value = -12414.667
$ws.Cells.Item(99, 14).Value = -7295.6665
$ws.Cells.Item(113, 8).Value = 989541.8
$ws.Cells.Item(113, 9).Value = 1088296
$ws.Cells.Item(113, 11).Value = 1088296
$ws.Cells.Item(113, 13).Value = -1086126
$ws.Cells.Item(126, 8).Value = 11509.417
$ws.Cells.Item(126, 9).Value = 13912.667
$ws.Cells.Item(126, 10).Value = 4299.6665
$ws.Cells.Item(126, 11).Value = 41738.001
$ws.Cells.Item(126, 12).Value = 12898.9995
$ws.Cells.Item(126, 13).Value = -39268.001
$ws.Cells.Item(126, 14).Value = -17838.9995
$ws.Cells.Item(132, 8).Value = 33337528
$ws.Cells.Item(132, 9).Value = 41670892
$ws.Cells.Item(132, 11).Value = 125012676
$ws.Cells.Item(132, 13).Value = -125010146
$ws.Cells.Item(134, 8).Value = 12503954
$ws.Cells.Item(134, 9).Value = 13892594
$ws.Cells.Item(134, 10).Value = 6199
$ws.Cells.Item(134, 11).Value = 41677782
$ws.Cells.Item(134, 12).Value = 18597
$ws.Cells.Item(134, 13).Value = -41675247
$ws.Cells.Item(134, 14).Value = -23667
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 3845.2236
$ws.Cells.Item(68, 9).Value = 1499.6666
$ws.Cells.Item(68, 10).Value = 3931.0366
$ws.Cells.Item(68, 11).Value = 4498.9998
$ws.Cells.Item(68, 12).Value = 11793.1098
$ws.Cells.Item(68, 13).Value = -3687.9998
$ws.Cells.Item(68, 14).Value = -13415.1098
$ws.Cells.Item(71, 8).Value = 3845.2236
$ws.Cells.Item(71, 9).Value = 1499.6666
$ws.Cells.Item(71, 10).Value = 3931.0366
$ws.Cells.Item(71, 11).Value = 13496.9994
$ws.Cells.Item(71, 12).Value = 35379.3294
$ws.Cells.Item(71, 13).Value = -9440.999400000001
$ws.Cells.Item(71, 14).Value = -43491.3294
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).ClearContents()
$ws.Cells.Item(92, 14).ClearContents()
$ws.Cells.Item(94, 8).Value = 12083.9375
$ws.Cells.Item(94, 10).Value = 13381.643
$ws.Cells.Item(94, 12).Value = 40144.929
$ws.Cells.Item(94, 14).Value = -41496.929
$ws.Cells.Item(113, 8).Value = 56306.055
$ws.Cells.Item(113, 10).Value = 875.7857
$ws.Cells.Item(113, 12).Value = 2627.3571
$ws.Cells.Item(113, 14).Value = -6967.3571
$ws.Cells.Item(117, 8).Value = 2986.8462
$ws.Cells.Item(117, 9).Value = 299.33334
$ws.Cells.Item(117, 10).Value = 3793.1
$ws.Cells.Item(117, 11).Value = 898.0000200000001
$ws.Cells.Item(117, 12).Value = 11379.3
$ws.Cells.Item(117, 13).Value = 2543.99998
$ws.Cells.Item(117, 14).Value = -18263.3
$ws.Cells.Item(120, 8).Value = 10000
$ws.Cells.Item(120, 9).Value = 10000
$ws.Cells.Item(120, 11).Value = 30000
$ws.Cells.Item(120, 13).Value = -25162
$ws.Cells.Item(132, 8).Value = 21733
$ws.Cells.Item(132, 9).Value = 21733
$ws.Cells.Item(132, 11).Value = 195597
$ws.Cells.Item(132, 13).Value = -193067
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(19, 8).Value = 12348.143
$ws.Cells.Item(19, 9).Value = 12348.143
$ws.Cells.Item(19, 11).Value = 12348.143
$ws.Cells.Item(19, 13).Value = -12060.143
$ws.Cells.Item(44, 8).Value = 39999
$ws.Cells.Item(44, 10).Value = 39999
$ws.Cells.Item(44, 12).Value = 39999
$ws.Cells.Item(44, 14).Value = -41191
$ws.Cells.Item(102, 8).Value = 3454.8572
$ws.Cells.Item(102, 9).Value = 3239.5789
$ws.Cells.Item(102, 11).Value = 3239.5789
$ws.Cells.Item(102, 13).Value = -1617.5789
$ws.Cells.Item(122, 8).Value = 4991412
$ws.Cells.Item(122, 9).Value = 5820396.5
$ws.Cells.Item(122, 10).Value = 17504
$ws.Cells.Item(122, 11).Value = 17461189.5
$ws.Cells.Item(122, 12).Value = 52512
$ws.Cells.Item(122, 13).Value = -17458739.5
$ws.Cells.Item(122, 14).Value = -57412
$ws.Cells.Item(126, 8).Value = 3013.5715
$ws.Cells.Item(126, 9).Value = 3013.5715
$ws.Cells.Item(126, 11).Value = 9040.7145
$ws.Cells.Item(126, 13).Value = -6570.7145
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 14719602
$ws.Cells.Item(132, 10).Value = 8799.200000000001
$ws.Cells.Item(132, 12).Value = 26397.6
$ws.Cells.Item(132, 14).Value = -31457.6
$ws.Cells.Item(141, 8).Value = 89997
$ws.Cells.Item(141, 10).Value = 89997
$ws.Cells.Item(141, 12).Value = 89997
$ws.Cells.Item(141, 14).Value = -100357
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2829.3845
$ws.Cells.Item(81, 10).Value = 2497
$ws.Cells.Item(81, 12).Value = 4994
$ws.Cells.Item(81, 14).Value = -7116
$ws.Cells.Item(84, 8).Value = 2829.3845
$ws.Cells.Item(84, 10).Value = 2497
$ws.Cells.Item(84, 12).Value = 24970
$ws.Cells.Item(84, 14).Value = -35578
$ws.Cells.Item(104, 8).Value = 20370
$ws.Cells.Item(104, 10).Value = 20370
$ws.Cells.Item(104, 12).Value = 20370
$ws.Cells.Item(104, 14).Value = -27358
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 14).ClearContents()
